$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.050.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.08%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.175.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.11%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.88%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.176.16"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.07%  "

$ws.Range("E9").Value = "  -0.77%  "

$ws.Range("E10").Value = "  -6.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.82%  "

$ws.Range("E12").Value = "  -3.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000237"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.53%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.700.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.06%  "

$ws.Range("E16").Value = "  -1.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.176.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.023.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.21%  "

$ws.Range("E19").Value = "  -4.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "460.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.36%  "

$ws.Range("E21").Value = "  -0.63%  "

$ws.Range("E22").Value = "  -6.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.77%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.43%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("E27").Value = "  +0.04%  "

$ws.Range("E28").Value = "  -4.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.64"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.15%  "

$ws.Range("E31").Value = "  -5.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.06%  "

$ws.Range("E33").Value = "  -4.08%  "

$ws.Range("E34").Value = "  -6.37%  "

$ws.Range("E35").Value = "  -6.97%  "

$ws.Range("E36").Value = "  -4.53%  "

$ws.Range("E37").Value = "  -2.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0707"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.40%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0387"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "404.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.16%  "

$ws.Range("E42").Value = "  -3.86%  "

$ws.Range("E43").Value = "  -6.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.814.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.45%  "

$ws.Range("E45").Value = "  -5.91%  "

$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("E47").Value = "  -5.93%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.95%  "

$ws.Range("E51").Value = "  -2.09%  "
